$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '327.90'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '-1.17%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '44.45'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '-0.47%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.365'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '-3.47%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.08363'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '1.20%'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '-4.97%'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9719'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '-0.53%'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '2.499'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '-5.47%'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.1121'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '-0.09%'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1900'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '-0.23%'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.09698'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '-3.19%'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.04604'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-1.40%'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.1061'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '0.22%'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.001290'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '0.48%'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.005894'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '-0.88%'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.361'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '-0.09%'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '4.414'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '-0.61%'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '8.543'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '-17.00%'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '0.15%'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '3.30%'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.04171'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '1.67%'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '-5.01%'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.004412'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '0.33%'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0002978'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '-20.48%'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '-2.45%'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05633'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '-2.02%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.007826'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '2.29%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.1412'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '-0.81%'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.007333'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '-3.02%'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '3.31%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.008712'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '4.78%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.3505'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00006907'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '-1.80%'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-0.16%'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '-7.98%'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.003530'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '39.85%'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00002100'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '-0.16%'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0002000'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '-0.16%'
